$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B11: change its string value from "R40" to "1" (kept as text,
# not converted to a number).
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
